$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Caulfield"
$ws.Range("B2").Value = "Metro Train - Frankston line"
$ws.Range("C2").Value = "30/12/20 4:30pm-17:00pm"
$ws.Range("D2").Value = "Case caught train from Caulfield to Cheltenham"
$ws.Range("E2").Value = "old"

$ws.Range("A3").Value = "Caulfield"
$ws.Range("B3").Value = "Metro Train - Frankston line"
$ws.Range("C3").Value = "30/12/20 4:30pm-5:00pm"
$ws.Range("D3").Value = "Case caught train from Caulfield to Cheltenham"
$ws.Range("E3").Value = "new"

$ws.Range("A:E").Select()
$ws.Range("A:E").EntireColumn.AutoFit() | Out-Null
